$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21; Excel shifts rows 21..46 down to 22..47
# (inheriting their existing values/formatting, including the D-column date style).
$ws.Rows(21).Insert()

# Populate the newly inserted row 21 with the new record.
$ws.Range("A21").Value = 11
$ws.Range("B21").Value = "Vega Monumental Concepción"
$ws.Range("C21").Value = "Bíobío"
$ws.Range("D21").Value = 44895
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = 100112026
$ws.Range("G21").Value = "Haba"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 14500
$ws.Range("N21").Value = "`$/saco 25 kilos"
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 580
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
